$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (for ownTeam and oppTeam),
# shifting the existing batsman..sr columns from D..I to F..K
$ws.Range("D1:E1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Data row
$ws.Range("D2").Value = "Rajasthan Royals"
$ws.Range("E2").Value = "Chennai Super Kings"
